# Update odds values on Sheet1 to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (La Equidad vs Santa Fe)
$ws.Range("G4").Value = 3.75
$ws.Range("I4").Value = 2.15
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.5
$ws.Range("U4").Value = 17
$ws.Range("V4").Value = 15
$ws.Range("AE4").Value = 9
$ws.Range("AF4").Value = 9.5
$ws.Range("AG4").Value = 19

# Row 5 (Macara vs Dep. Cuenca)
$ws.Range("J5").Value = 1.11
$ws.Range("K5").Value = 6.5

# Row 6 (Manta vs Tecnico U.)
$ws.Range("I6").Value = 2.05
$ws.Range("K6").Value = 9
$ws.Range("X6").Value = 29
$ws.Range("AE6").Value = 9.5
$ws.Range("AG6").Value = 19
$ws.Range("AH6").Value = 19

# Row 15 (Club Tijuana vs Atl. San Luis)
$ws.Range("G15").Value = 1.75
$ws.Range("H15").Value = 4.75
$ws.Range("I15").Value = 3.5
$ws.Range("R15").Value = 1.4
$ws.Range("S15").Value = 2.75

# Row 19 (Sportivo Trinidense vs Atl. Tembetary)
$ws.Range("G19").Value = 2.2
$ws.Range("I19").Value = 3.3
$ws.Range("J19").Value = 1.08
$ws.Range("K19").Value = 8
$ws.Range("N19").Value = 2.2
$ws.Range("O19").Value = 1.65
$ws.Range("T19").Value = 7
$ws.Range("U19").Value = 10
$ws.Range("W19").Value = 21
$ws.Range("X19").Value = 21
$ws.Range("Z19").Value = 7.5
$ws.Range("AD19").Value = 9
$ws.Range("AE19").Value = 15
$ws.Range("AG19").Value = 34
$ws.Range("AH19").Value = 29

# Row 29 (Monagas vs Rayo Zuliano)
$ws.Range("G29").Value = 1.65
$ws.Range("H29").Value = 3.5
$ws.Range("I29").Value = 5
$ws.Range("M29").Value = 2.82
$ws.Range("N29").Value = 1.93
$ws.Range("O29").Value = 1.7
$ws.Range("R29").Value = 1.9
$ws.Range("S29").Value = 1.72
$ws.Range("T29").Value = 6.2
$ws.Range("U29").Value = 7.4
$ws.Range("W29").Value = 12.5
$ws.Range("X29").Value = 14
$ws.Range("AA29").Value = 6.8
$ws.Range("AB29").Value = 17.5
$ws.Range("AE29").Value = 28
$ws.Range("AF29").Value = 16.5
$ws.Range("AG29").Value = 90
$ws.Range("AH29").Value = 55
$ws.Range("AI29").Value = 65
$ws.Range("AJ29").Value = 900
